$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 31
$ws.Range("I2").Value = 90
$ws.Range("J2").Value = 376
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 116
$ws.Range("M2").Value = 3
$ws.Range("N2").Value = 66
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 7
$ws.Range("S2").Value = 36
$ws.Range("T2").Value = 73
$ws.Range("U2").Value = 4
$ws.Range("V2").Value = 622
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 592
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 12
$ws.Range("AA2").Value = 5
